$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source "Price" column (D) is plain text in the original workbook
# (coinranking.com formats with dots as thousands separators, e.g.
# "54.413.32"), and the "Volume(1h)" column (E) always carries leading/
# trailing padding spaces around the percentage. Writing such strings
# straight into Range.Value would normally be fine, but plain decimal-
# looking values (e.g. "22.60") get silently re-interpreted by Excel as
# numbers, which would both change the cell type and drop meaningful
# trailing zeros (22.60 -> 22.6). Prefixing with a leading apostrophe is
# the standard Excel convention to force literal text entry so the value
# is stored exactly as written.

$ws.Range("D2").Value = "54.413.32"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "2.278.07"
$ws.Range("E3").Value = "  +2.26%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'498.59"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("D6").Value = "'128.26"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("E9").Value = "  +2.65%  "

$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("E11").Value = "  +3.29%  "

$ws.Range("D12").Value = "'4.71"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").Value = "2.683.49"
$ws.Range("E13").Value = "  +2.24%  "

$ws.Range("D14").Value = "'22.60"
$ws.Range("E14").Value = "  +5.50%  "

$ws.Range("D15").Value = "54.304.89"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").Value = "2.272.71"
$ws.Range("E17").Value = "  +1.79%  "

$ws.Range("D18").Value = "'10.28"
$ws.Range("E18").Value = "  +5.19%  "

$ws.Range("D19").Value = "'4.13"
$ws.Range("E19").Value = "  +2.12%  "

$ws.Range("D20").Value = "'304.54"
$ws.Range("E20").Value = "  +2.52%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").Value = "'61.85"
$ws.Range("E23").Value = "  -2.60%  "

$ws.Range("D24").Value = "'0.995"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("E25").Value = "  +2.67%  "

$ws.Range("E26").Value = "  +2.89%  "

$ws.Range("D27").Value = "'176.20"
$ws.Range("E27").Value = "  +8.33%  "

$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").Value = "'5.96"
$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("E31").Value = "  +1.38%  "

$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").Value = "'17.76"
$ws.Range("E33").Value = "  +2.08%  "

$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'0.922"
$ws.Range("E35").Value = "  +9.51%  "

$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("D37").Value = "'3.75"
$ws.Range("E37").Value = "  +3.41%  "

$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("E40").Value = "  +1.57%  "

$ws.Range("D41").Value = "'125.68"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "'4.77"
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("D43").Value = "'0.0491"
$ws.Range("E43").Value = "  +2.74%  "

$ws.Range("D44").Value = "'0.0898"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").Value = "'240.03"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("E47").Value = "  -0.55%  "

$ws.Range("D48").Value = "'0.0206"
$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("E49").Value = "  +1.03%  "

$ws.Range("D50").Value = "'16.28"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("E51").Value = "  +0.33%  "
